# Add a new global-constant entry "MaxGuideQuestId" with a temporary
# value of 4 as the next row of the GlobalConstantIntTable sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "MaxGuideQuestId"
$ws.Cells.Item($newRow, 2).Value = 4
